$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 755.44446
$ws.Cells.Item(8, 9).Value = 828.4286
$ws.Cells.Item(8, 11).Value = 2485.2858
$ws.Cells.Item(8, 13).Value = -2346.2858

$ws.Cells.Item(19, 8).Value = 913.3
$ws.Cells.Item(19, 10).Value = 764.6667
$ws.Cells.Item(19, 12).Value = 764.6667
$ws.Cells.Item(19, 14).Value = -1114.6667

$ws.Cells.Item(20, 8).Value = 1662.25
$ws.Cells.Item(20, 9).Value = 1662.25
$ws.Cells.Item(20, 11).Value = 1662.25
$ws.Cells.Item(20, 13).Value = -1432.25

$ws.Cells.Item(35, 8).Value = 1662.25
$ws.Cells.Item(35, 9).Value = 1662.25
$ws.Cells.Item(35, 11).Value = 1662.25
$ws.Cells.Item(35, 13).Value = -1283.25

$ws.Cells.Item(70, 8).Value = 137503760
$ws.Cells.Item(70, 9).Value = 33334166
$ws.Cells.Item(70, 11).Value = 100002498
$ws.Cells.Item(70, 13).Value = -100002228

$ws.Cells.Item(73, 8).Value = 137503760
$ws.Cells.Item(73, 9).Value = 33334166
$ws.Cells.Item(73, 11).Value = 100002498
$ws.Cells.Item(73, 13).Value = -100001562

$ws.Cells.Item(129, 8).Value = 1565.9375
$ws.Cells.Item(129, 9).Value = 1058.1818
$ws.Cells.Item(129, 11).Value = 3174.5454
$ws.Cells.Item(129, 13).Value = 1825.4546

$ws.Cells.Item(138, 8).Value = 1950.9803
$ws.Cells.Item(138, 10).Value = 2118.0557
$ws.Cells.Item(138, 12).Value = 6354.1671
$ws.Cells.Item(138, 14).Value = -16634.1671

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4213.087
$ws.Cells.Item(32, 9).Value = 3498.2046
$ws.Cells.Item(32, 11).Value = 3498.2046
$ws.Cells.Item(32, 13).Value = -3211.2046

$ws.Cells.Item(45, 8).Value = 2068.4443
$ws.Cells.Item(45, 9).Value = 1325.6
$ws.Cells.Item(45, 11).Value = 1325.6
$ws.Cells.Item(45, 13).Value = -948.5999999999999

$ws.Cells.Item(74, 8).Value = 2552.8823
$ws.Cells.Item(74, 9).Value = 2877.8462
$ws.Cells.Item(74, 10).Value = 1496.75
$ws.Cells.Item(74, 11).Value = 2877.8462
$ws.Cells.Item(74, 12).Value = 1496.75
$ws.Cells.Item(74, 13).Value = -2003.8462
$ws.Cells.Item(74, 14).Value = -3244.75

$ws.Cells.Item(77, 8).Value = 2552.8823
$ws.Cells.Item(77, 9).Value = 2877.8462
$ws.Cells.Item(77, 10).Value = 1496.75
$ws.Cells.Item(77, 11).Value = 14389.231
$ws.Cells.Item(77, 12).Value = 7483.75
$ws.Cells.Item(77, 13).Value = -10021.231
$ws.Cells.Item(77, 14).Value = -16219.75

$ws.Cells.Item(102, 8).Value = 3722.7878
$ws.Cells.Item(102, 9).Value = 2629.875
$ws.Cells.Item(102, 11).Value = 2629.875
$ws.Cells.Item(102, 13).Value = -1007.875

$ws.Cells.Item(132, 8).Value = 8309.25
$ws.Cells.Item(132, 9).Value = 9317.786
$ws.Cells.Item(132, 10).Value = 1249.5
$ws.Cells.Item(132, 11).Value = 27953.358
$ws.Cells.Item(132, 12).Value = 3748.5
$ws.Cells.Item(132, 13).Value = -25423.358
$ws.Cells.Item(132, 14).Value = -8808.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 822.5
$ws.Cells.Item(22, 9).Value = 893.8
$ws.Cells.Item(22, 10).Value = 466
$ws.Cells.Item(22, 11).Value = 893.8
$ws.Cells.Item(22, 12).Value = 466
$ws.Cells.Item(22, 13).Value = -720.8
$ws.Cells.Item(22, 14).Value = -812

$ws.Cells.Item(45, 8).Value = 25600
$ws.Cells.Item(45, 10).Value = 25600
$ws.Cells.Item(45, 12).Value = 25600
$ws.Cells.Item(45, 14).Value = -27216

$ws.Cells.Item(99, 8).Value = 2218.7307
$ws.Cells.Item(99, 9).Value = 1245.8667
$ws.Cells.Item(99, 10).Value = 3545.3635
$ws.Cells.Item(99, 11).Value = 1245.8667
$ws.Cells.Item(99, 12).Value = 3545.3635
$ws.Cells.Item(99, 13).Value = 252.1333
$ws.Cells.Item(99, 14).Value = -6541.363499999999

$ws.Cells.Item(134, 8).Value = 2808.9473
$ws.Cells.Item(134, 9).Value = 2945.9167
$ws.Cells.Item(134, 11).Value = 8837.750100000001
$ws.Cells.Item(134, 13).Value = -6302.750100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2935.1316
$ws.Cells.Item(31, 9).Value = 1021.6
$ws.Cells.Item(31, 11).Value = 1021.6
$ws.Cells.Item(31, 13).Value = -726.6

$ws.Cells.Item(34, 8).Value = 2935.1316
$ws.Cells.Item(34, 9).Value = 1021.6
$ws.Cells.Item(34, 11).Value = 1021.6
$ws.Cells.Item(34, 13).Value = -819.6

$ws.Cells.Item(99, 8).Value = 2710.125
$ws.Cells.Item(99, 9).Value = 2754.4285
$ws.Cells.Item(99, 11).Value = 2754.4285
$ws.Cells.Item(99, 13).Value = -1256.4285

$ws.Cells.Item(126, 8).Value = 2710.125
$ws.Cells.Item(126, 9).Value = 2754.4285
$ws.Cells.Item(126, 11).Value = 8263.2855
$ws.Cells.Item(126, 13).Value = -5793.2855

$ws.Cells.Item(132, 8).Value = 1795.6538
$ws.Cells.Item(132, 9).Value = 1651.9048
$ws.Cells.Item(132, 11).Value = 4955.7144
$ws.Cells.Item(132, 13).Value = -2425.7144

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 133.5
$ws.Cells.Item(33, 10).Value = 179.77777
$ws.Cells.Item(33, 12).Value = 1078.66662
$ws.Cells.Item(33, 14).Value = -1644.66662

$ws.Cells.Item(44, 8).Value = 2077
$ws.Cells.Item(44, 10).Value = 1425
$ws.Cells.Item(44, 12).Value = 4275
$ws.Cells.Item(44, 14).Value = -5071

$ws.Cells.Item(47, 8).Value = 243.44444
$ws.Cells.Item(47, 9).Value = 263.2
$ws.Cells.Item(47, 10).Value = 218.75
$ws.Cells.Item(47, 11).Value = 789.5999999999999
$ws.Cells.Item(47, 12).Value = 656.25
$ws.Cells.Item(47, 13).Value = -358.5999999999999
$ws.Cells.Item(47, 14).Value = -1518.25

$ws.Cells.Item(51, 8).Value = 1383.3334
$ws.Cells.Item(51, 9).Value = 800
$ws.Cells.Item(51, 10).Value = 1500
$ws.Cells.Item(51, 11).Value = 2400
$ws.Cells.Item(51, 12).Value = 4500
$ws.Cells.Item(51, 13).Value = -1940
$ws.Cells.Item(51, 14).Value = -5420

$ws.Cells.Item(121, 8).Value = 665.4666999999999
$ws.Cells.Item(121, 9).Value = 253
$ws.Cells.Item(121, 10).Value = 1026.375
$ws.Cells.Item(121, 11).Value = 759
$ws.Cells.Item(121, 12).Value = 3079.125
$ws.Cells.Item(121, 13).Value = 551
$ws.Cells.Item(121, 14).Value = -5699.125

$ws.Cells.Item(132, 8).Value = 2858.0833
$ws.Cells.Item(132, 10).Value = 2888.889
$ws.Cells.Item(132, 12).Value = 26000.001
$ws.Cells.Item(132, 14).Value = -31060.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 12).ClearContents()
$ws.Cells.Item(19, 14).Value = 0

$ws.Cells.Item(132, 8).Value = 3146
$ws.Cells.Item(132, 9).Value = 3012.5
$ws.Cells.Item(132, 10).Value = 3199.4
$ws.Cells.Item(132, 11).Value = 9037.5
$ws.Cells.Item(132, 12).Value = 9598.200000000001
$ws.Cells.Item(132, 13).Value = -6507.5
$ws.Cells.Item(132, 14).Value = -14658.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3252.2104
$ws.Cells.Item(46, 9).Value = 3166.75
$ws.Cells.Item(46, 11).Value = 3166.75
$ws.Cells.Item(46, 13).Value = -2978.75

$ws.Cells.Item(87, 8).Value = 45000
$ws.Cells.Item(87, 10).Value = 45000
$ws.Cells.Item(87, 12).Value = 45000
$ws.Cells.Item(87, 14).Value = -47246

$ws.Cells.Item(90, 8).Value = 45000
$ws.Cells.Item(90, 10).Value = 45000
$ws.Cells.Item(90, 12).Value = 135000
$ws.Cells.Item(90, 14).Value = -146232

$ws.Cells.Item(128, 8).Value = 50143
$ws.Cells.Item(128, 10).Value = 50143
$ws.Cells.Item(128, 12).Value = 50143
$ws.Cells.Item(128, 14).Value = -60103

$ws.Cells.Item(132, 8).Value = 4467.7144
$ws.Cells.Item(132, 9).Value = 4149.875
$ws.Cells.Item(132, 11).Value = 12449.625
$ws.Cells.Item(132, 13).Value = -9919.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 16336.333
$ws.Cells.Item(41, 9).Value = 13884.5
$ws.Cells.Item(41, 10).Value = 18297.8
$ws.Cells.Item(41, 11).Value = 13884.5
$ws.Cells.Item(41, 12).Value = 18297.8
$ws.Cells.Item(41, 13).Value = -13494.5
$ws.Cells.Item(41, 14).Value = -19077.8

$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).ClearContents()
$ws.Cells.Item(75, 13).ClearContents()
$ws.Cells.Item(75, 14).Value = 0

$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).ClearContents()
$ws.Cells.Item(78, 13).ClearContents()
$ws.Cells.Item(78, 14).Value = 0

$ws.Cells.Item(122, 8).Value = 4113.68
$ws.Cells.Item(122, 9).Value = 2352.7144
$ws.Cells.Item(122, 10).Value = 6354.909
$ws.Cells.Item(122, 11).Value = 7058.1432
$ws.Cells.Item(122, 12).Value = 19064.727
$ws.Cells.Item(122, 13).Value = -4608.1432
$ws.Cells.Item(122, 14).Value = -23964.727

$ws.Cells.Item(132, 8).Value = 4960.5615
$ws.Cells.Item(132, 9).Value = 5829.4224
$ws.Cells.Item(132, 11).Value = 17488.2672
$ws.Cells.Item(132, 13).Value = -14958.2672
